$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.677.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.51%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.680.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "523.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.578"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.699.98"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.43"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.07%  "
$ws.Range("E11").Value = "  -0.89%  "
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("E13").Value = "  +2.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.152.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.610.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.49%  "
$ws.Range("E16").Value = "  +0.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000138"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.703.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "350.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.25%  "
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.55%  "
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("E26").Value = "  +4.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.997"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0815"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.85%  "
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("E33").Value = "  +0.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "148.62"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.26"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.51%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.61%  "
$ws.Range("B37").Value = "SuiNetwork"
$ws.Range("C37").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.947"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.879"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.51"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.91"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.42%  "
$ws.Range("E41").Value = "  -2.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "282.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.03"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.61%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0990"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.609"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.138.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.43%  "
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0539"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.90"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.37%  "
$ws.Range("E50").Value = "  +0.83%  "
$ws.Range("E51").Value = "  +1.77%  "
